$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temp staging cell used to write percent-like text without Excel
# auto-converting "NN%" strings into numeric percentages (which would
# change the cell type/style). We format the staging cell as Text,
# copy it, paste-special (values only) into the destination, then
# clear the staging cell so it leaves no trace in the sheet.
$stage = $ws.Range("Z1")

$ws.Range("E2").Value = "2026-02-21 22:48:16"
$stage.NumberFormat = "@"
$stage.Value = "45%"
$stage.Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("E3").Value = "2026-02-21 22:48:19"
$ws.Range("E4").Value = "2026-02-21 22:48:21"
$stage.NumberFormat = "@"
$stage.Value = "74%"
$stage.Copy()
$ws.Range("H4").PasteSpecial(-4163)
$ws.Range("O4").Value = "9.1 °C"
$ws.Range("E5").Value = "2026-02-21 22:48:23"
$ws.Range("E6").Value = "2026-02-21 22:48:26"
$ws.Range("E7").Value = "2026-02-21 22:48:28"
$ws.Range("E8").Value = "2026-02-21 22:48:31"
$ws.Range("O8").Value = "11.3 °C"
$ws.Range("E9").Value = "2026-02-21 22:48:33"
$stage.NumberFormat = "@"
$stage.Value = "59%"
$stage.Copy()
$ws.Range("H9").PasteSpecial(-4163)
$ws.Range("K9").Value = "15.0 MJ/m2"
$ws.Range("N9").Value = "5.4 °C 22:28 TU"
$ws.Range("O9").Value = "12.9 °C"
$ws.Range("E10").Value = "2026-02-21 22:48:34"
$ws.Range("O10").Value = "8.3 °C"
$ws.Range("E11").Value = "2026-02-21 22:48:35"
$ws.Range("O11").Value = "8.6 °C"
$ws.Range("E12").Value = "2026-02-21 22:48:36"
$stage.NumberFormat = "@"
$stage.Value = "65%"
$stage.Copy()
$ws.Range("H12").PasteSpecial(-4163)
$ws.Range("N12").Value = "6.6 °C 22:23 TU"
$ws.Range("O12").Value = "12.3 °C"
$ws.Range("E13").Value = "2026-02-21 22:48:38"
$ws.Range("J13").Value = "1032.0 hPa"
$ws.Range("E14").Value = "2026-02-21 22:48:39"
$stage.NumberFormat = "@"
$stage.Value = "72%"
$stage.Copy()
$ws.Range("H14").PasteSpecial(-4163)
$ws.Range("O14").Value = "11.1 °C"
$ws.Range("E15").Value = "2026-02-21 22:48:40"
$stage.NumberFormat = "@"
$stage.Value = "58%"
$stage.Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("N15").Value = "4.9 °C 22:00 TU"
$ws.Range("O15").Value = "12.6 °C"
$ws.Range("E16").Value = "2026-02-21 22:48:42"
$ws.Range("E17").Value = "2026-02-21 22:48:45"
$ws.Range("K17").Value = "16.8 MJ/m2"
$ws.Range("E18").Value = "2026-02-21 22:48:47"
$stage.NumberFormat = "@"
$stage.Value = "76%"
$stage.Copy()
$ws.Range("H18").PasteSpecial(-4163)
$ws.Range("O18").Value = "8.4 °C"
$ws.Range("E19").Value = "2026-02-21 22:48:49"
$ws.Range("E20").Value = "2026-02-21 22:48:51"
$ws.Range("K20").Value = "16.2 MJ/m2"
$ws.Range("O20").Value = "2.9 °C"
$ws.Range("E21").Value = "2026-02-21 22:48:52"
$ws.Range("E22").Value = "2026-02-21 22:48:54"
$ws.Range("G22").Value = "110 cm"
$ws.Range("E23").Value = "2026-02-21 22:48:56"
$ws.Range("E24").Value = "2026-02-21 22:48:59"
$ws.Range("O24").Value = "6.2 °C"
$ws.Range("E25").Value = "2026-02-21 22:49:01"
$ws.Range("E26").Value = "2026-02-21 22:49:04"
$ws.Range("J26").Value = "1027.6 hPa"
$ws.Range("E27").Value = "2026-02-21 22:49:06"
$ws.Range("E28").Value = "2026-02-21 22:49:09"
$stage.NumberFormat = "@"
$stage.Value = "74%"
$stage.Copy()
$ws.Range("H28").PasteSpecial(-4163)
$ws.Range("J28").Value = "1029.8 hPa"
$ws.Range("O28").Value = "8.0 °C"
$ws.Range("E29").Value = "2026-02-21 22:49:11"
$stage.NumberFormat = "@"
$stage.Value = "69%"
$stage.Copy()
$ws.Range("H29").PasteSpecial(-4163)
$ws.Range("N29").Value = "5.1 °C 22:17 TU"
$ws.Range("O29").Value = "11.2 °C"
$ws.Range("E30").Value = "2026-02-21 22:49:14"
$stage.NumberFormat = "@"
$stage.Value = "69%"
$stage.Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("O30").Value = "11.3 °C"
$ws.Range("E31").Value = "2026-02-21 22:49:16"
$ws.Range("O31").Value = "12.4 °C"
$ws.Range("E32").Value = "2026-02-21 22:49:19"
$ws.Range("O32").Value = "4.6 °C"
$ws.Range("E33").Value = "2026-02-21 22:49:21"
$stage.NumberFormat = "@"
$stage.Value = "53%"
$stage.Copy()
$ws.Range("H33").PasteSpecial(-4163)
$ws.Range("O33").Value = "6.5 °C"
$ws.Range("E34").Value = "2026-02-21 22:49:23"
$stage.NumberFormat = "@"
$stage.Value = "40%"
$stage.Copy()
$ws.Range("H34").PasteSpecial(-4163)
$ws.Range("O34").Value = "4.2 °C"
$ws.Range("E35").Value = "2026-02-21 22:49:26"
$ws.Range("E36").Value = "2026-02-21 22:49:28"
$stage.NumberFormat = "@"
$stage.Value = "60%"
$stage.Copy()
$ws.Range("H36").PasteSpecial(-4163)
$ws.Range("O36").Value = "12.9 °C"
$ws.Range("E37").Value = "2026-02-21 22:49:31"
$stage.NumberFormat = "@"
$stage.Value = "75%"
$stage.Copy()
$ws.Range("H37").PasteSpecial(-4163)
$ws.Range("O37").Value = "5.5 °C"
$ws.Range("E38").Value = "2026-02-21 22:49:33"
$ws.Range("O38").Value = "9.6 °C"
$ws.Range("E39").Value = "2026-02-21 22:49:36"
$stage.NumberFormat = "@"
$stage.Value = "32%"
$stage.Copy()
$ws.Range("H39").PasteSpecial(-4163)
$ws.Range("E40").Value = "2026-02-21 22:49:38"
$stage.NumberFormat = "@"
$stage.Value = "55%"
$stage.Copy()
$ws.Range("H40").PasteSpecial(-4163)
$ws.Range("O40").Value = "8.3 °C"
$ws.Range("E41").Value = "2026-02-21 22:49:40"
$ws.Range("O41").Value = "11.1 °C"
$ws.Range("E42").Value = "2026-02-21 22:49:43"
$ws.Range("E43").Value = "2026-02-21 22:49:45"
$stage.NumberFormat = "@"
$stage.Value = "76%"
$stage.Copy()
$ws.Range("H43").PasteSpecial(-4163)
$ws.Range("O43").Value = "7.1 °C"
$ws.Range("E44").Value = "2026-02-21 22:49:48"
$ws.Range("K44").Value = "15.9 MJ/m2"
$ws.Range("O44").Value = "2.1 °C"
$ws.Range("E45").Value = "2026-02-21 22:49:50"
$stage.NumberFormat = "@"
$stage.Value = "67%"
$stage.Copy()
$ws.Range("H45").PasteSpecial(-4163)
$ws.Range("J45").Value = "1032.4 hPa"
$ws.Range("E46").Value = "2026-02-21 22:49:52"
$ws.Range("O46").Value = "9.5 °C"

$stage.Clear()
$excel.CutCopyMode = $false

